$d = $word.ActiveDocument

# 1. Remove the hidden "_GoBack" bookmark left over from the previous edit session.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Append new content at the end of the body (after the existing trailing
#    empty paragraph, before the sectPr): an "AND/OR" paragraph, an empty
#    paragraph, an "Elementer mangler innhold." paragraph, and a final empty
#    paragraph.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = "<w:p $wNs/>" +
       "<w:p $wNs><w:r><w:t>AND/OR</w:t></w:r></w:p>" +
       "<w:p $wNs/>" +
       "<w:p $wNs><w:r><w:t>Elementer mangler innhold.</w:t></w:r></w:p>" +
       "<w:p $wNs/>"

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(1)
$r.InsertXML($xml)
